# Rename the first sheet ("opus_big Validation") and make it the active/
# selected tab with a new active-cell selection, matching the edits
# captured in the commit: the workbook's active tab moves from the last
# sheet ("opus_big LSP Fine aWCE ") to the first sheet, which is also
# renamed to "opus_big Pure WCE + Baselines".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Name = "opus_big Pure WCE + Baselines"

# Activate this sheet so it becomes the selected/visible tab (this also
# clears the tabSelected flag that was previously on the last sheet).
$ws.Select()

# Move the in-sheet selection/active cell to D45, matching the saved
# selection in the workbook.
$ws.Range("D45").Select()
